$d = $word.ActiveDocument

$replacements = @(
    @("601×5=", "140×6="),
    @("756×8=", "539×4="),
    @("693×4=", "201×3="),
    @("884×8=", "118×8="),
    @("958×3=", "478×3="),
    @("158×7=", "444×7="),
    @("830×3=", "283×9="),
    @("686×9=", "959×8="),
    @("602×3=", "565×6="),
    @("738×9=", "861×3="),
    @("187×9=", "299×5="),
    @("378×9=", "421×3="),
    @("265×3=", "414×3="),
    @("685×7=", "288×3="),
    @("388×7=", "988×2="),
    @("810×9=", "577×6="),
    @("885×5=", "723×2="),
    @("190×6=", "266×2="),
    @("149×6=", "220×2="),
    @("124×8=", "828×6="),
    @("150×3=", "981×9="),
    @("461×2=", "610×4="),
    @("833×4=", "411×8="),
    @("389×7=", "530×9="),
    @("610×7=", "109×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
